# feat: add 2022-Q4 data
#
# Inserts a new "2022-Q4" worksheet (positioned right after the "总计"
# summary sheet, before the existing "2022-Q3" sheet) with the new
# quarter's fund-holding figures, and inserts a matching new row at the
# top of the "总计" sheet's data table.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "总计" (totals) sheet: insert a new row 2 for 2022-Q4, pushing the
#    existing quarters down by one row.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

$total.Rows.Item(2).Insert()
$total.Range("B2:D2").ClearFormats()

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 4.59

# Renumber the index column (A) for the rows that shifted down, and give
# A2 the same style as the rest of the index column.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A3").Copy($total.Range("A2"))
$total.Range("A2").Value = 0

# ---------------------------------------------------------------------
# 2) Add the new "2022-Q4" worksheet right after "总计", cloning the
#    layout/styles of the existing "2022-Q3" sheet and then overwriting
#    the data with the new quarter's figures.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)
$newSheet = $wb.Worksheets.Add($null, $total)
$newSheet.Name = "2022-Q4"

$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Range("A1:H3").Copy($newSheet.Range("A1:H3"))

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'012348"
$newSheet.Range("C2").Value = "天弘恒生科技指数（QDII）A"
$newSheet.Range("D2").Value = "'39.65"
$newSheet.Range("E2").Value = "'93.67"
$newSheet.Range("F2").Value = "'5.95"
$newSheet.Range("G2").Value = "'2.3592"
$newSheet.Range("H2").Value = 6

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "'012349"
$newSheet.Range("C3").Value = "天弘恒生科技指数（QDII）C"
$newSheet.Range("D3").Value = "'37.52"
$newSheet.Range("E3").Value = "'93.67"
$newSheet.Range("F3").Value = "'5.95"
$newSheet.Range("G3").Value = "'2.2324"
$newSheet.Range("H3").Value = 6
